$d = $word.ActiveDocument

# Update the date heading
$d.Paragraphs.Item(1).Range.Text = '2025-03-12 Wednesday'

# Update the table cells (20 rows x 5 columns)
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = '37+47='
$t.Cell(1, 2).Range.Text = '87-19='
$t.Cell(1, 3).Range.Text = '20-3='
$t.Cell(1, 4).Range.Text = '91-26='
$t.Cell(1, 5).Range.Text = '33-5='
$t.Cell(2, 1).Range.Text = '92-23='
$t.Cell(2, 2).Range.Text = '37+15='
$t.Cell(2, 3).Range.Text = '29+59='
$t.Cell(2, 4).Range.Text = '45-26='
$t.Cell(2, 5).Range.Text = '76-19='
$t.Cell(3, 1).Range.Text = '23-4='
$t.Cell(3, 2).Range.Text = '92-28='
$t.Cell(3, 3).Range.Text = '34+19='
$t.Cell(3, 4).Range.Text = '80-73='
$t.Cell(3, 5).Range.Text = '55-7='
$t.Cell(4, 1).Range.Text = '82-19='
$t.Cell(4, 2).Range.Text = '48+35='
$t.Cell(4, 3).Range.Text = '52-33='
$t.Cell(4, 4).Range.Text = '26+9='
$t.Cell(4, 5).Range.Text = '36+47='
$t.Cell(5, 1).Range.Text = '64-25='
$t.Cell(5, 2).Range.Text = '13+8='
$t.Cell(5, 3).Range.Text = '94-19='
$t.Cell(5, 4).Range.Text = '86-37='
$t.Cell(5, 5).Range.Text = '80-34='
$t.Cell(6, 1).Range.Text = '26+5='
$t.Cell(6, 2).Range.Text = '24+27='
$t.Cell(6, 3).Range.Text = '44-5='
$t.Cell(6, 4).Range.Text = '42-23='
$t.Cell(6, 5).Range.Text = '90-78='
$t.Cell(7, 1).Range.Text = '58+25='
$t.Cell(7, 2).Range.Text = '7+89='
$t.Cell(7, 3).Range.Text = '41-14='
$t.Cell(7, 4).Range.Text = '15+58='
$t.Cell(7, 5).Range.Text = '68+3='
$t.Cell(8, 1).Range.Text = '90-31='
$t.Cell(8, 2).Range.Text = '19+57='
$t.Cell(8, 3).Range.Text = '70-22='
$t.Cell(8, 4).Range.Text = '94-25='
$t.Cell(8, 5).Range.Text = '29+6='
$t.Cell(9, 1).Range.Text = '32-9='
$t.Cell(9, 2).Range.Text = '64-17='
$t.Cell(9, 3).Range.Text = '65-57='
$t.Cell(9, 4).Range.Text = '38+55='
$t.Cell(9, 5).Range.Text = '32+49='
$t.Cell(10, 1).Range.Text = '60-58='
$t.Cell(10, 2).Range.Text = '33+29='
$t.Cell(10, 3).Range.Text = '48+27='
$t.Cell(10, 4).Range.Text = '28+53='
$t.Cell(10, 5).Range.Text = '5+79='
$t.Cell(11, 1).Range.Text = '11-9='
$t.Cell(11, 2).Range.Text = '17+74='
$t.Cell(11, 3).Range.Text = '70-19='
$t.Cell(11, 4).Range.Text = '7+59='
$t.Cell(11, 5).Range.Text = '19+57='
$t.Cell(12, 1).Range.Text = '18+75='
$t.Cell(12, 2).Range.Text = '93-46='
$t.Cell(12, 3).Range.Text = '23+9='
$t.Cell(12, 4).Range.Text = '72-48='
$t.Cell(12, 5).Range.Text = '53-18='
$t.Cell(13, 1).Range.Text = '81-13='
$t.Cell(13, 2).Range.Text = '9+18='
$t.Cell(13, 3).Range.Text = '14+68='
$t.Cell(13, 4).Range.Text = '74+19='
$t.Cell(13, 5).Range.Text = '30-2='
$t.Cell(14, 1).Range.Text = '14+49='
$t.Cell(14, 2).Range.Text = '46-7='
$t.Cell(14, 3).Range.Text = '31-9='
$t.Cell(14, 4).Range.Text = '68+7='
$t.Cell(14, 5).Range.Text = '93-58='
$t.Cell(15, 1).Range.Text = '80-26='
$t.Cell(15, 2).Range.Text = '38+39='
$t.Cell(15, 3).Range.Text = '46+8='
$t.Cell(15, 4).Range.Text = '52-45='
$t.Cell(15, 5).Range.Text = '84-68='
$t.Cell(16, 1).Range.Text = '71-33='
$t.Cell(16, 2).Range.Text = '63-9='
$t.Cell(16, 3).Range.Text = '29+64='
$t.Cell(16, 4).Range.Text = '60-26='
$t.Cell(16, 5).Range.Text = '92-78='
$t.Cell(17, 1).Range.Text = '8+43='
$t.Cell(17, 2).Range.Text = '70-11='
$t.Cell(17, 3).Range.Text = '61-13='
$t.Cell(17, 4).Range.Text = '92-69='
$t.Cell(17, 5).Range.Text = '33-6='
$t.Cell(18, 1).Range.Text = '8+39='
$t.Cell(18, 2).Range.Text = '21-2='
$t.Cell(18, 3).Range.Text = '90-26='
$t.Cell(18, 4).Range.Text = '8+3='
$t.Cell(18, 5).Range.Text = '67+25='
$t.Cell(19, 1).Range.Text = '85-28='
$t.Cell(19, 2).Range.Text = '9+48='
$t.Cell(19, 3).Range.Text = '5+57='
$t.Cell(19, 4).Range.Text = '34-19='
$t.Cell(19, 5).Range.Text = '54-19='
$t.Cell(20, 1).Range.Text = '72-9='
$t.Cell(20, 2).Range.Text = '3+8='
$t.Cell(20, 3).Range.Text = '41-35='
$t.Cell(20, 4).Range.Text = '3+78='
$t.Cell(20, 5).Range.Text = '81-13='
